$d = $word.ActiveDocument

# --- Step 1: split the <dryingUrl>...</dryingUrl> paragraph's single run into
# three runs, turning the middle (url-ish) portion into a real hyperlink. ---
$p3 = $d.Paragraphs.Item(3).Range
$p3Start = $p3.Start
$linkRange = $d.Range($p3Start + 11, $p3Start + 122)
$d.Hyperlinks.Add($linkRange, "https://www.wikihow.com/Dry-Figs#:~:text=To%20dry%20figs%2C%20start%20by,for%20up%20to%2036%20hours.") | Out-Null

# --- Step 2: insert a new empty paragraph followed by a paragraph containing
# the imgur picture link, right after the dryingUrl paragraph. ---
$p4 = $d.Paragraphs.Item(4).Range
$insertPoint = $d.Range($p4.Start, $p4.Start)
$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p/><w:p><w:r><w:t>https://i.imgur.com/qwa4ybm.jpg</w:t></w:r></w:p><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$insertPoint.InsertXML($xml)
# InsertXML's trailing empty paragraph merges-in as a spare - remove it.
$extra = $d.Paragraphs.Item(6).Range
$extra.Delete()

# --- Step 3: register the "Hyperlink" / "Unresolved Mention" character
# styles that Word mints the first time a hyperlink is inserted. ---
$hs = $d.Styles.Add("Hyperlink", 2)
$hs.BaseStyle = "DefaultParagraphFont"
$hs.Priority = 99
$hs.UnhideWhenUsed = $true
$hs.Font.TextColor.ObjectThemeColor = 10
$hs.Font.Underline = 1

$ums = $d.Styles.Add("UnresolvedMention", 2)
$ums.NameLocal = "Unresolved Mention"
$ums.BaseStyle = "DefaultParagraphFont"
$ums.Priority = 99
$ums.UnhideWhenUsed = $true
$ums.Font.Color = 0x5C5E60
